$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.825.45'
$ws.Range("E2").Value = '  +2.89%  '
$ws.Range("D3").Value = '3.439.82'
$ws.Range("E3").Value = '  +2.18%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''579.91'
$ws.Range("E5").Value = '  +4.50%  '
$ws.Range("D6").Value = '''188.42'
$ws.Range("E6").Value = '  +8.37%  '
$ws.Range("D7").Value = '''0.630'
$ws.Range("E7").Value = '  -0.37%  '
$ws.Range("D8").Value = '3.436.15'
$ws.Range("E8").Value = '  +2.44%  '
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  -1.99%  '
$ws.Range("D11").Value = '''0.643'
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").Value = '''56.78'
$ws.Range("E12").Value = '  +5.98%  '
$ws.Range("E13").Value = '  -1.60%  '
$ws.Range("E14").Value = '  +2.95%  '
$ws.Range("D15").Value = '3.992.19'
$ws.Range("E15").Value = '  +2.29%  '
$ws.Range("D16").Value = '''18.70'
$ws.Range("E16").Value = '  +2.70%  '
$ws.Range("D17").Value = '3.444.08'
$ws.Range("E17").Value = '  +2.41%  '
$ws.Range("D18").Value = '66.852.16'
$ws.Range("E18").Value = '  +2.91%  '
$ws.Range("D19").Value = '''12.05'
$ws.Range("E19").Value = '  +1.56%  '
$ws.Range("E20").Value = '  -1.70%  '
$ws.Range("E21").Value = '  +2.63%  '
$ws.Range("E22").Value = '  +6.96%  '
$ws.Range("E23").Value = '  +7.81%  '
$ws.Range("D24").Value = '''16.84'
$ws.Range("E24").Value = '  +23.27%  '
$ws.Range("D25").Value = '''4.34'
$ws.Range("E25").Value = '  +6.80%  '
$ws.Range("D26").Value = '''89.15'
$ws.Range("E26").Value = '  +2.29%  '
$ws.Range("E27").Value = '  +3.12%  '
$ws.Range("D28").Value = '''10.96'
$ws.Range("E28").Value = '  +2.25%  '
$ws.Range("D29").Value = '''9.03'
$ws.Range("E29").Value = '  +4.36%  '
$ws.Range("D30").Value = '''31.18'
$ws.Range("E30").Value = '  +0.19%  '
$ws.Range("D31").Value = '''7.35'
$ws.Range("E31").Value = '  +12.50%  '
$ws.Range("D32").Value = '''11.75'
$ws.Range("E32").Value = '  +2.82%  '
$ws.Range("D33").Value = '''596.73'
$ws.Range("E33").Value = '  +3.64%  '
$ws.Range("D34").Value = '''64.04'
$ws.Range("E34").Value = '  +1.87%  '
$ws.Range("E35").Value = '  +4.04%  '
$ws.Range("E36").Value = '  +5.67%  '
$ws.Range("D37").Value = '''1.00'
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").Value = '''36.81'
$ws.Range("E38").Value = '  +3.38%  '
$ws.Range("D39").Value = '''0.385'
$ws.Range("E39").Value = '  +4.07%  '
$ws.Range("E40").Value = '  -2.96%  '
$ws.Range("D41").Value = '0.0₃0755'
$ws.Range("E41").Value = '  +2.04%  '
$ws.Range("D42").Value = '3.236.79'
$ws.Range("E42").Value = '  +5.03%  '
$ws.Range("E43").Value = '  +5.08%  '
$ws.Range("D44").Value = '''0.0430'
$ws.Range("E44").Value = '  +3.29%  '
$ws.Range("D45").Value = '''2.86'
$ws.Range("E45").Value = '  +26.86%  '
$ws.Range("E46").Value = '  +3.78%  '
$ws.Range("D47").Value = '''3.38'
$ws.Range("E47").Value = '  +17.72%  '
$ws.Range("D48").Value = '''3.20'
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("E49").Value = '  +0.29%  '
$ws.Range("D50").Value = '''1.00'
$ws.Range("E50").Value = '  +0.11%  '
$ws.Range("D51").Value = '''8.60'
$ws.Range("E51").Value = '  +4.07%  '
